$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("__footings__")

# Clear old columns H:K (end_point/column_name/stable were removed, old row_start..col_end columns shift left)
$ws.Range("H1:K9").ClearContents()

# Row 1
$ws.Cells.Item(1, 1).Value = 'worksheet'
$ws.Cells.Item(1, 2).Value = 'mapping'
$ws.Cells.Item(1, 3).Value = 'dtype'
$ws.Cells.Item(1, 4).Value = 'row_start'
$ws.Cells.Item(1, 5).Value = 'col_start'
$ws.Cells.Item(1, 6).Value = 'row_end'
$ws.Cells.Item(1, 7).Value = 'col_end'

# Row 2
$ws.Cells.Item(2, 1).Value = 'test_dict'
$ws.Cells.Item(2, 2).Value = '/outer'
$ws.Cells.Item(2, 3).Value = '<class ''str''>'
$ws.Cells.Item(2, 4).Value = 2
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 2
$ws.Cells.Item(2, 7).Value = 2

# Row 3
$ws.Cells.Item(3, 1).Value = 'test_dict'
$ws.Cells.Item(3, 2).Value = '/outer/inner'
$ws.Cells.Item(3, 3).Value = '<class ''str''>'
$ws.Cells.Item(3, 4).Value = 2
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 2
$ws.Cells.Item(3, 7).Value = 3

# Row 4
$ws.Cells.Item(4, 1).Value = 'test_dict'
$ws.Cells.Item(4, 2).Value = '/outer/inner/endpoint1'
$ws.Cells.Item(4, 3).Value = '<class ''str''>'
$ws.Cells.Item(4, 4).Value = 2
$ws.Cells.Item(4, 5).Value = 4
$ws.Cells.Item(4, 6).Value = 2
$ws.Cells.Item(4, 7).Value = 4

# Row 5
$ws.Cells.Item(5, 1).Value = 'test_dict'
$ws.Cells.Item(5, 2).Value = '/outer/inner/endpoint1'
$ws.Cells.Item(5, 3).Value = '<class ''int''>'
$ws.Cells.Item(5, 4).Value = 2
$ws.Cells.Item(5, 5).Value = 5
$ws.Cells.Item(5, 6).Value = 2
$ws.Cells.Item(5, 7).Value = 5

# Row 6
$ws.Cells.Item(6, 1).Value = 'test_dict'
$ws.Cells.Item(6, 2).Value = '/outer/inner/endpoint2'
$ws.Cells.Item(6, 3).Value = '<class ''str''>'
$ws.Cells.Item(6, 4).Value = 3
$ws.Cells.Item(6, 5).Value = 4
$ws.Cells.Item(6, 6).Value = 3
$ws.Cells.Item(6, 7).Value = 4

# Row 7
$ws.Cells.Item(7, 1).Value = 'test_dict'
$ws.Cells.Item(7, 2).Value = '/outer/inner/endpoint2'
$ws.Cells.Item(7, 3).Value = '<class ''int''>'
$ws.Cells.Item(7, 4).Value = 3
$ws.Cells.Item(7, 5).Value = 5
$ws.Cells.Item(7, 6).Value = 3
$ws.Cells.Item(7, 7).Value = 5

# Row 8
$ws.Cells.Item(8, 1).Value = 'test_dict'
$ws.Cells.Item(8, 2).Value = '/endpoint3'
$ws.Cells.Item(8, 3).Value = '<class ''str''>'
$ws.Cells.Item(8, 4).Value = 4
$ws.Cells.Item(8, 5).Value = 2
$ws.Cells.Item(8, 6).Value = 4
$ws.Cells.Item(8, 7).Value = 2

# Row 9
$ws.Cells.Item(9, 1).Value = 'test_dict'
$ws.Cells.Item(9, 2).Value = '/endpoint3'
$ws.Cells.Item(9, 3).Value = '<class ''int''>'
$ws.Cells.Item(9, 4).Value = 4
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 4
$ws.Cells.Item(9, 7).Value = 3

